# Generate Report for Handback
# The file "606488d1-f3c8-41ce-9f66-3363d3eb6950.md" has now been handed back
# (in sync with en-US) for both the zh-cn and de-de locales. Update the
# Overview sheet and the per-locale report sheets accordingly.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# ---- Overview sheet: row 2 is the 606488d1 file ----
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("E2").Value = $status
$ovw.Range("F2").Value = $status

# ---- zh-cn sheet: row 2 is the 606488d1 file ----
$zh = $wb.Worksheets.Item("zh-cn")
$zh.Range("C2").Value = $status
$zh.Range("I2").Value = "606488d1-f3c8-41ce-9f66-3363d3eb6950.md"
$zh.Range("J2").Value = "606488d1-f3c8-41ce-9f66-3363d3eb6950.bf3dca8ee1989928cb4cb1ce9f4606c33618c2fb.zh-cn.xlf"
$zh.Range("K2").Value = "2016-09-06 17:14:25"
$zh.Hyperlinks.Add($zh.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80d88d3b910e630e5e6c9dc30008b9ec5951c85b/e2e/606488d1-f3c8-41ce-9f66-3363d3eb6950.md", "", "", "606488d1-f3c8-41ce-9f66-3363d3eb6950.md")

# ---- de-de sheet: row 2 is the 606488d1 file ----
$de = $wb.Worksheets.Item("de-de")
$de.Range("C2").Value = $status
$de.Range("I2").Value = "606488d1-f3c8-41ce-9f66-3363d3eb6950.md"
$de.Range("J2").Value = "606488d1-f3c8-41ce-9f66-3363d3eb6950.bf3dca8ee1989928cb4cb1ce9f4606c33618c2fb.de-de.xlf"
$de.Range("K2").Value = "2016-09-06 17:14:33"
$de.Hyperlinks.Add($de.Range("I2"), "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/80d88d3b910e630e5e6c9dc30008b9ec5951c85b/e2e/606488d1-f3c8-41ce-9f66-3363d3eb6950.md", "", "", "606488d1-f3c8-41ce-9f66-3363d3eb6950.md")

# ---- Column widths: widen the columns that now hold longer text ----
$ovw.Columns.Item(5).ColumnWidth = 29.9777047293527
$ovw.Columns.Item(6).ColumnWidth = 29.9777047293527

$zh.Columns.Item(3).ColumnWidth = 29.9777047293527
$zh.Columns.Item(9).ColumnWidth = 40
$zh.Columns.Item(10).ColumnWidth = 40

$de.Columns.Item(3).ColumnWidth = 29.9777047293527
$de.Columns.Item(9).ColumnWidth = 40
$de.Columns.Item(10).ColumnWidth = 40
